# Generate Report for Handoff
# Updates the localization-status workbook with the new handoff file id
# (3d929e4c-fc2f-420a-9dc6-836075e1a685 -> 940c231c-cdc6-458a-8955-ab89043e0388)
# and refreshed handoff timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$newId = "940c231c-cdc6-458a-8955-ab89043e0388"
$newZhHash = "989096b01a9e7fe8fde66d198afa168a9d196fcb"

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newId.md"
}
$wsOverview.Range("G2").Value = "2016-08-14 17:16:50"

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newId.md"
foreach ($h in $wsZh.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}
$wsZh.Range("G2").Value = "$newId.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-14 17:16:41"

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newId.md"
foreach ($h in $wsDe.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}
$wsDe.Range("G2").Value = "$newId.$newZhHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-14 17:16:50"
